$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1931
$ws1.Range("F3").Value = 528
$ws1.Range("F6").Value = 2876
$ws1.Range("F10").Value = 1600
$ws1.Range("F11").Value = 1581
$ws1.Range("F25").Value = 277
$ws1.Range("F26").Value = 89
$ws1.Range("G26").Value = 36.6
$ws1.Range("F28").Value = 5
$ws1.Range("F29").Value = 1871
$ws1.Range("F33").Value = 124
$ws1.Range("F34").Value = 574

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1931
$ws4.Range("F4").Value = 528
$ws4.Range("F7").Value = 2876
$ws4.Range("F11").Value = 1600
$ws4.Range("F12").Value = 1581
$ws4.Range("F26").Value = 277
$ws4.Range("F27").Value = 90
$ws4.Range("G27").Value = 36.6
$ws4.Range("F29").Value = 5
$ws4.Range("F30").Value = 1871
$ws4.Range("F34").Value = 124
$ws4.Range("F35").Value = 574
